# Bug Fixes and Minor Updates
# - proxy=0 / proxy=n condition fix reflected by changing F2 threshold value from 20 to 0
# - Company name input (A2) reset to a fresh example "yes bank"
# - The old example company list that used to live in column A (A3:A16) is moved out
#   to its own "temp" helper sheet so it no longer clutters the input form.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

# ---------------------------------------------------------------------------
# 1. Create the new "temp" sheet (placed right after "input") and move the
#    old example company names (previously in input!A2:A16) onto it.
# ---------------------------------------------------------------------------
$oldCompanies = @(
    "tata steel bsl",
    "tata motors",
    "bhushan steel",
    "bharti airtel",
    "indiabulls ventures",
    "crisil",
    "bank of baroda",
    "wipro",
    "bharat petroleum corp",
    "icici bank",
    "xyz123",
    "atulya123",
    "coal india ltd",
    "dewan housing",
    "pidilite industries"
)

$tempSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$tempSheet.Name = "temp"

for ($i = 0; $i -lt $oldCompanies.Length; $i++) {
    $tempSheet.Cells.Item($i + 1, 1).Value = $oldCompanies[$i]
}

# Match the look of the old input column (width/fill) as closely as possible.
$tempSheet.Columns.Item(1).ColumnWidth = 20.25
$tempSheet.Range("A1:A16").Interior.Color = $ws.Range("A2").Interior.Color

# Reproduce the saved selection (A1:A16 selected) on the temp sheet, then
# return focus back to the input sheet so it stays the active tab.
[void]$tempSheet.Range("A1:A16").Select()

# ---------------------------------------------------------------------------
# 2. Update the "input" sheet itself.
# ---------------------------------------------------------------------------

# Clear out the old example company names from A3:A16 (the cells themselves
# are removed, not just blanked).
[void]$ws.Range("A3:A16").Clear()

# A2 now holds a fresh example company name.
$ws.Range("A2").Value = "yes bank"

# The proxy-related threshold in F2 goes from 20 down to 0.
$ws.Range("F2").Value = 0

# Restore "input" as the active sheet/tab and update the saved cursor
# position + zoom level to match the edited workbook.
[void]$ws.Select()
[void]$ws.Range("F5").Select()
$excel.ActiveWindow.Zoom = 85
